# Update the date heading paragraph.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-07-03 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-07-04 Friday", 2)

# Update each table cell by (row, column) position so that values which
# collide with other cells' old/new text do not get clobbered by a
# whole-document Find/Replace.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "28÷6=4, 4"
$t.Cell(1,2).Range.Text  = "93÷6=15, 3"
$t.Cell(1,3).Range.Text  = "19÷4=4, 3"
$t.Cell(1,4).Range.Text  = "80÷3=26, 2"
$t.Cell(1,5).Range.Text  = "53÷5=10, 3"

$t.Cell(5,1).Range.Text  = "82÷8=10, 2"
$t.Cell(5,2).Range.Text  = "52÷7=7, 3"
$t.Cell(5,3).Range.Text  = "70÷9=7, 7"
$t.Cell(5,4).Range.Text  = "64÷4=16, 0"
$t.Cell(5,5).Range.Text  = "41÷9=4, 5"

$t.Cell(9,1).Range.Text  = "40÷7=5, 5"
$t.Cell(9,2).Range.Text  = "86÷3=28, 2"
$t.Cell(9,3).Range.Text  = "89÷5=17, 4"
$t.Cell(9,4).Range.Text  = "41÷9=4, 5"
$t.Cell(9,5).Range.Text  = "57÷9=6, 3"

$t.Cell(13,1).Range.Text = "51÷6=8, 3"
$t.Cell(13,2).Range.Text = "12÷7=1, 5"
$t.Cell(13,3).Range.Text = "89÷2=44, 1"
$t.Cell(13,4).Range.Text = "77÷8=9, 5"
$t.Cell(13,5).Range.Text = "35÷6=5, 5"

$t.Cell(17,1).Range.Text = "75÷9=8, 3"
$t.Cell(17,2).Range.Text = "93÷7=13, 2"
$t.Cell(17,3).Range.Text = "23÷5=4, 3"
$t.Cell(17,4).Range.Text = "39÷6=6, 3"
$t.Cell(17,5).Range.Text = "72÷8=9, 0"
